{"js": "// Update each three-digit-by-one-digit multiplication prompt in place.\n// Each old value is unique in the document, so a plain text search +\n// \"Replace\" swaps exactly one table cell per pair.\nconst replacements = [\n  [\"652\u00d74=\", \"753\u00d73=\"],\n  [\"860\u00d78=\", \"347\u00d74=\"],\n  [\"694\u00d74=\", \"723\u00d78=\"],\n  [\"220\u00d78=\", \"588\u00d78=\"],\n  [\"864\u00d75=\", \"833\u00d79=\"],\n  [\"511\u00d78=\", \"160\u00d72=\"],\n  [\"290\u00d74=\", \"250\u00d79=\"],\n  [\"147\u00d77=\", \"866\u00d79=\"],\n  [\"300\u00d78=\", \"630\u00d78=\"],\n  [\"523\u00d74=\", \"881\u00d74=\"],\n  [\"446\u00d77=\", \"512\u00d79=\"],\n  [\"963\u00d78=\", \"108\u00d76=\"],\n  [\"775\u00d77=\", \"709\u00d78=\"],\n  [\"299\u00d77=\", \"427\u00d72=\"],\n  [\"107\u00d78=\", \"221\u00d76=\"],\n  [\"285\u00d77=\", \"299\u00d78=\"],\n  [\"434\u00d79=\", \"701\u00d73=\"],\n  [\"262\u00d72=\", \"260\u00d72=\"],\n  [\"451\u00d72=\", \"546\u00d74=\"],\n  [\"283\u00d74=\", \"508\u00d78=\"],\n  [\"425\u00d76=\", \"212\u00d76=\"],\n  [\"578\u00d73=\", \"617\u00d73=\"],\n  [\"404\u00d76=\", \"538\u00d73=\"],\n  [\"576\u00d74=\", \"740\u00d72=\"],\n  [\"729\u00d75=\", \"805\u00d77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update each three-digit-by-one-digit multiplication prompt in place.\n# Each old value is unique in the document, so a bounded Find/Replace\n# (MatchCase + MatchWholeWord, no wildcards) swaps exactly one cell each.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"652\u00d74=\", $true, $true, $false, $false, $false, $true, 1, $false, \"753\u00d73=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"860\u00d78=\", $true, $true, $false, $false, $false, $true, 1, $false, \"347\u00d74=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"694\u00d74=\", $true, $true, $false, $false, $false, $true, 1, $false, \"723\u00d78=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"220\u00d78=\", $true, $true, $false, $false, $false, $true, 1, $false, \"588\u00d78=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"864\u00d75=\", $true, $true, $false, $false, $false, $true, 1, $false, \"833\u00d79=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"511\u00d78=\", $true, $true, $false, $false, $false, $true, 1, $false, \"160\u00d72=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"290\u00d74=\", $true, $true, $false, $false, $false, $true, 1, $false, \"250\u00d79=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"147\u00d77=\", $true, $true, $false, $false, $false, $true, 1, $false, \"866\u00d79=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"300\u00d78=\", $true, $true, $false, $false, $false, $true, 1, $false, \"630\u00d78=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"523\u00d74=\", $true, $true, $false, $false, $false, $true, 1, $false, \"881\u00d74=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"446\u00d77=\", $true, $true, $false, $false, $false, $true, 1, $false, \"512\u00d79=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"963\u00d78=\", $true, $true, $false, $false, $false, $true, 1, $false, \"108\u00d76=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"775\u00d77=\", $true, $true, $false, $false, $false, $true, 1, $false, \"709\u00d78=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"299\u00d77=\", $true, $true, $false, $false, $false, $true, 1, $false, \"427\u00d72=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"107\u00d78=\", $true, $true, $false, $false, $false, $true, 1, $false, \"221\u00d76=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"285\u00d77=\", $true, $true, $false, $false, $false, $true, 1, $false, \"299\u00d78=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"434\u00d79=\", $true, $true, $false, $false, $false, $true, 1, $false, \"701\u00d73=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"262\u00d72=\", $true, $true, $false, $false, $false, $true, 1, $false, \"260\u00d72=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"451\u00d72=\", $true, $true, $false, $false, $false, $true, 1, $false, \"546\u00d74=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"283\u00d74=\", $true, $true, $false, $false, $false, $true, 1, $false, \"508\u00d78=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"425\u00d76=\", $true, $true, $false, $false, $false, $true, 1, $false, \"212\u00d76=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"578\u00d73=\", $true, $true, $false, $false, $false, $true, 1, $false, \"617\u00d73=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"404\u00d76=\", $true, $true, $false, $false, $false, $true, 1, $false, \"538\u00d73=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"576\u00d74=\", $true, $true, $false, $false, $false, $true, 1, $false, \"740\u00d72=\", 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"729\u00d75=\", $true, $true, $false, $false, $false, $true, 1, $false, \"805\u00d77=\", 2) | Out-Null\n\n"}
